$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 was a blank placeholder record (A10="9", B10/C10 stored as the
# text "0", the rest empty). Row 11 held the "牟礼いこいの家" record
# (A11="10") and row 12 held the "牟礼いこいの広場" record (A12="11").
# The real edit drops the placeholder row's B:N content and pulls the
# two real records up one row each (each row keeps its own column-A
# sequence number), so the trailing now-duplicate row (12) goes away and
# the sheet shrinks from A1:N12 to A1:N11.

# --- Row 10 <- old Row 11 (B:M), column A ("9") is untouched ---
# B/C hold digit-only text ("34.34630528" / "134.12677778"); format the
# cell as Text first so Excel keeps them as strings instead of coercing
# them to numbers, then drop the format override so no stray number
# format is left behind on the cell.
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "34.34630528"
$ws.Range("B10").ClearFormats()

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "134.12677778"
$ws.Range("C10").ClearFormats()

$ws.Range("D10").Value = "老人いこいの家牟礼いこいの家"
$ws.Range("E10").Value = "高松市牟礼町牟礼1978-1"
$ws.Range("J10").Value = "月火水木金土日"
$ws.Range("K10").Value = "09:00"
$ws.Range("L10").Value = "18:00"
$ws.Range("M10").Value = "時間外利用可能"

# --- Row 11 <- old Row 12 (B:E); J:M clear out to blank ---
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "34.33712806"
$ws.Range("B11").ClearFormats()

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "134.13068722"
$ws.Range("C11").ClearFormats()

$ws.Range("D11").Value = "老人いこいの家牟礼いこいの広場"
$ws.Range("E11").Value = "高松市牟礼町牟礼568"
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""

# The old row 12 is now a duplicate trailing row; remove it so the sheet
# ends at row 11 (A1:N11).
$ws.Rows("12").Delete()
